$wb = $excel.ActiveWorkbook

# The "as_of_utc" column (AA) on both data sheets gets refreshed to the
# latest publish timestamp for every data row (rows 2-26).
$timestamp = "2025-11-20 07:05:41"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("AA2:AA26").Value = $timestamp
}
